# Update the "Duration of contains()" chart data (Chart sheet!B2:B25) with
# new measurements, and relabel the value axis to call out the unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:B25 (row 2 = SNM 8 ... row 25 = SNM 31)
$newValues = @(
    1746222.9,
    852565.1,
    254353.3,
    110736.2,
    57767.7,
    31158.9,
    15566.6,
    8341.4,
    4252.6,
    2525.6,
    1087.1,
    538.1,
    281.7,
    152.1,
    80.3,
    42.7,
    24.6,
    15.4,
    21.3,
    11.7,
    10.9,
    7.0,
    7.6,
    6.6
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Keep the chart's series in sync with the refreshed worksheet data.
$chartObject = $ws.ChartObjects().Item(1)
$chart = $chartObject.Chart
$series = $chart.SeriesCollection(1)
$series.Values = $ws.Range("B2:B25")

# Rename the value axis title to clarify the unit of measurement.
$chart.Axes(2).AxisTitle.Text = "Duration in Microseconds"
